$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 ("R40") is replaced with the text value "1".
# Assigning the literal "1" directly would be auto-coerced to a number by
# Excel, so we build it as a text formula first and then collapse the
# formula down to a static value (copy / paste-special values), which keeps
# the cell's text ("string") type and its existing style.
$cell = $ws.Range("B11")
$cell.Formula = "=""1"""
$cell.Copy()
$cell.PasteSpecial(-4163)
